# 5.1.1. EQUIPOS BAJO TRÁMITE RMA — ajuste de formatos
#
# The paragraph that used to hold the freeform "Bajo el trámite de RMA
# (Return Merchandise Authorization) ... no se tramita equipos bajo el
# proceso de garantía" sentence is wiped out and reformatted as an empty
# "Heading 1" (pStyle Ttulo1) paragraph: tab stop at 1410 twips (70.5pt),
# SpaceBefore 165 twips (8.25pt), FirstLineIndent reset to 0.

$d = $word.ActiveDocument

# Locate the target paragraph robustly by its distinctive trailing text
# rather than a hard-coded index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*no se tramita equipos bajo el proceso de garant*") {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    # Strip the paragraph's runs, keeping the paragraph mark itself.
    $r = $target.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = ""

    # Re-style the (now empty) paragraph as a Heading 1 / Ttulo1 paragraph.
    $target.Range.Style = $d.Styles.Item("Heading 1")

    # Direct paragraph formatting overrides matching the new pPr.
    $target.Format.SpaceBefore = 8.25
    $target.Format.FirstLineIndent = 0

    # Explicit custom tab stop at 1410 twips (70.5 pt).
    $target.Range.Select()
    $word.Selection.ParagraphFormat.TabStops.Add(70.5)

    Write-Output "Reformatted RMA paragraph."
} else {
    Write-Output "Target RMA paragraph not found."
}
